$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are stored as text (coin prices use locale-specific
# formatting such as "27.884.80" that Excel would otherwise reinterpret
# as a number), so force a text format before assigning each new value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.884.80"
$ws.Range("E2").Value = "  +1.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.752.82"
$ws.Range("E3").Value = "  -0.84%  "

$ws.Range("E4").Value = "  -0.54%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.45"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("E6").Value = "  -0.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3817"
$ws.Range("E7").Value = "  -1.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3389"
$ws.Range("E8").Value = "  -0.94%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.66"
$ws.Range("E9").Value = "  -5.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.111"
$ws.Range("E10").Value = "  -3.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07207"
$ws.Range("E11").Value = "  -3.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.41"
$ws.Range("E13").Value = "  -0.73%  "

$ws.Range("E14").Value = "  -3.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.131"
$ws.Range("E15").Value = "  +0.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.756.35"
$ws.Range("E16").Value = "  -1.13%  "

$ws.Range("E17").Value = "  -2.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06613"
$ws.Range("E18").Value = "  -1.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "78.91"
$ws.Range("E19").Value = "  -4.61%  "

$ws.Range("E20").Value = "  -0.65%  "

$ws.Range("E21").Value = "  -4.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.213"
$ws.Range("E22").Value = "  -3.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.903.04"
$ws.Range("E23").Value = "  +1.16%  "

$ws.Range("E24").Value = "  -4.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.384"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.09"
$ws.Range("E26").Value = "  -1.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.78"
$ws.Range("E27").Value = "  -5.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.310"
$ws.Range("E28").Value = "  -5.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.958.13"
$ws.Range("E29").Value = "  -0.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.271"
$ws.Range("E30").Value = "  -11.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.89"
$ws.Range("E31").Value = "  -2.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.022"
$ws.Range("E32").Value = "  +1.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.814"
$ws.Range("E33").Value = "  -5.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08796"
$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.18"
$ws.Range("E35").Value = "  -4.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6570"
$ws.Range("E36").Value = "  -3.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06163"
$ws.Range("E37").Value = "  -3.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02282"
$ws.Range("E38").Value = "  -6.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.118"
$ws.Range("E39").Value = "  -5.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.507"
$ws.Range("E40").Value = "  -2.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2100"
$ws.Range("E41").Value = "  -4.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.205"
$ws.Range("E42").Value = "  -3.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.964"
$ws.Range("E43").Value = "  -5.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  -0.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.70"
$ws.Range("E45").Value = "  -3.33%  "

$ws.Range("E46").Value = "  -0.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6044"
$ws.Range("E47").Value = "  -3.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.05"
$ws.Range("E48").Value = "  -4.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.996"
$ws.Range("E49").Value = "  -5.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.121"
$ws.Range("E50").Value = "  +5.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.169"
$ws.Range("E51").Value = "  +1.88%  "
